$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('3535 Opal Meadow Heights Aged Care Community Meadow Heights', 27),
    @('95 Napier Street Apartment Complex Fitzroy', 5),
    @('Adorn Cosmetics Clayton', 5),
    @('Al Haj Halal Meats Glenroy', 27),
    @('Al-Taqwa College Truganina', 9),
    @('Amiga Montessori Craigieburn', 26),
    @('Australia Post Distribution Centre Sunshine West', 5),
    @('Baxter Foods Australia Campbellfield', 5),
    @('Best&Less Fountain Gate Narre Warren', 5),
    @('Budget Car and Truck Rentals Campbellfield', 6),
    @('CS Square Caroline Springs', 11),
    @('Campbellfield Ford Complex Vaccination Clinic Campbellfield', 14),
    @('Cannie Road Construction Site Cannie', 5),
    @('Cannie Road Construction Site Cannie', 6),
    @('Caroline Springs Police Station', 11),
    @('Cedars Medical Clinic Coburg', 34),
    @('Chemist Warehouse Campbellfield DC', 6),
    @('Chemist Warehouse Fillo Drive Somerton', 11),
    @('City of Moreland Community', 5),
    @('City of Wyndham Community', 6),
    @('Classy Cabinets and Kitchens Craigieburn', 9),
    @('Coles Aurora Village Epping', 6),
    @('Coles Broadmeadows Central Shopping Centre', 7),
    @('Coles Campbellfield Plaza Campbellfield', 12),
    @('Coles Coburg North Village', 31),
    @('Coles Greenvale Shopping Centre', 6),
    @('Coles Pakenham Place Shopping Centre', 13),
    @('Coles Roxburgh Village Roxburgh Park', 9),
    @('Community Kids Bayswater Early Education Centre Bayswater North', 11),
    @('Community Kids Meadow Heights', 11),
    @('Construction Site Olea Apartment Caulfield North', 17),
    @('Costco Wholesale Epping', 31),
    @('Crusader Caravans Epping', 23),
    @('DayHab Rehabilitation Treatment Centre Ringwood East', 5),
    @('Direct Freight Express Cambellfield', 13),
    @('Don Watson Coldstore Derrimut', 5),
    @('Epworth Healthcare Epworth Richmond Emergency Department', 6),
    @('Fine Food Holdings Pty Ltd Dandenong South', 8),
    @('Fitzroy Community School Fitzroy North', 33),
    @('Fonterra Manufacturing Workplace Campbellfield', 10),
    @('General Foods Campbellfield', 13),
    @('Glenroy West Primary School', 6),
    @('Goodstart Early Learning Altona', 11),
    @('Green Leaves Early Learning Cairnlea', 5),
    @('Green Leaves Early Learning Centre Highlands Craigieburn', 11),
    @('Gumboots Early Learning Centre South Morang', 5),
    @('Hamilton Marino 236 Jasper Road McKinnon', 13),
    @('Health Care Providers Association South Melbourne', 7),
    @('Hello Fresh Warehouse Ravenhall', 5),
    @('IGA Meadow Heights Shopping Centre Meadow Heights', 6),
    @('ISS Factory Level 1 Terminal 2 Melbourne Airport Tullamarine', 11),
    @('Ibis Kingsgate Hotel Melbourne', 6),
    @('Ilim College Glenroy Campus Hadfield', 10),
    @('Ilim Learning Sanctuary Glenroy', 9),
    @('Industrial Galvanizers Valmont Coatings Campbellfield', 22),
    @('KFC Fawkner', 6),
    @('Kasr Sweets Coolaroo', 7),
    @('Kids House Early Learning Cheltenham', 9),
    @('Kippers Seafood Werribee', 6),
    @('Kool Kidz Childcare Narre Warren', 6),
    @('Learning Nest Early Learning Centre Meadow Heights', 5),
    @('Level Crossing Removal Project Lilydale Construction Site John Street', 7),
    @('Lineage Logistics Laverton North', 8),
    @('Linfox Somerton National Distribution Centre Somerton', 9),
    @('McDonald''s Craigieburn North', 5),
    @('Mecca D.C Warehouse Melbourne Airport', 9),
    @('Melbourne Assessment Prison West Melbourne', 5),
    @('Melbourne Metropolitan Remand Centre Ravenhall', 11),
    @('Melbourne Truck Repairs Campbellfield', 7),
    @('Melbourne West Police Station Docklands', 7),
    @('Mercy Hospital for Women Heidelberg', 5),
    @('Mernda YMCA Early Learning Centre Mernda', 5),
    @('Mill Park Police Station Mill Park', 7),
    @('MyCentre Childcare Broadmeadows', 17),
    @('National Gallery of Victoria Melbourne', 9),
    @('Nido Early School Ascot Vale', 6),
    @('Nido Early School Glenroy', 11),
    @('Nido Early School Moonee Ponds', 13),
    @('Northern Health Northern Hospital Epping Emergency Department Tier 1B', 55),
    @('Northern Health The Northern Hospital Epping', 8),
    @('OnQ Plumbing and Excavations Craigieburn', 19),
    @('Oporto Coolaroo', 11),
    @('Oscar Romero Catholic Primary School Craigieburn', 6),
    @('Our Lady Help of Christian''s Primary School Brunswick East', 10),
    @('Paisley Park Early Learning Centre Bundoora', 6),
    @('Panorama Construction Site Whitehorse Rd Box Hill', 13),
    @('Private Residence Northern Community Services Fawkner', 5),
    @('Ramsay Health Care Warringal Private Hospital Heidelberg', 9),
    @('Richmond Quarter 261-271 Bridge Road Construction Site Richmond', 11),
    @('Sacca''s Fruit World Broadmeadows Central Shopping Centre', 6),
    @('Salta Drive Construction Site Rangedale Drainage Altona North', 7),
    @('St Margaret''s Primary School OSHC Maribyrnong', 11),
    @('St Vincents Hospital Emergency Department Melbourne', 9),
    @('Sultan Halal Meats & Poultry Campbellfield', 5),
    @('Tek Foods Somerton', 20),
    @('The Homestead Child and Family Centre Roxburgh Park', 10),
    @('The Huntly-Goornong Rail Works', 5),
    @('The Royal Children''s Hospital Melbourne Emergency Department Parkville Tier 1B', 18),
    @('The Royal Melbourne Hospital AMU Ward Parkville', 27),
    @('ThorwestenCabinets Pakenham', 14),
    @('Total Window Concepts Hoppers Crossing', 6),
    @('Unilodge College Square Student Accommodation 570 Lygon Street Carlton', 12),
    @('Wallaby Childcare Wollert', 11),
    @('Werribee Mercy Hospital Emergency Department', 11),
    @('Western Health Footscray Hospital Emergency Department', 6),
    @('Western Health Sunshine Hospital Emergency Department', 9),
    @('Woodlands Long Day Care and Kindergarten Roxburgh Park', 5),
    @('Woolworths Greenvale Lakes Roxburgh Park', 5),
    @('Yara Childcare Centre Truganina', 8)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
